{"js": "// Update the date line and every \"A\u00d7B=\" expression in the practice table.\n// We walk the document's paragraphs in document order and replace the\n// text of each paragraph that currently holds one of the \"before\" values\n// with its corresponding \"after\" value. Matching on the exact current\n// text (rather than a blind positional index) keeps this script correct\n// even if the document's paragraph layout differs slightly from what we\n// expect, while still being driven entirely by paragraph identity rather\n// than string search-and-replace (so there is no risk of a later\n// replacement accidentally re-matching text inserted by an earlier one).\nconst replacements = new Map([\n  [\"2025-08-26 Tuesday\", \"2025-08-27 Wednesday\"],\n  [\"644\u00d75=\", \"925\u00d73=\"],\n  [\"234\u00d76=\", \"101\u00d78=\"],\n  [\"284\u00d72=\", \"705\u00d78=\"],\n  [\"766\u00d74=\", \"205\u00d79=\"],\n  [\"437\u00d78=\", \"687\u00d75=\"],\n  [\"404\u00d73=\", \"610\u00d74=\"],\n  [\"463\u00d73=\", \"279\u00d79=\"],\n  [\"396\u00d79=\", \"160\u00d72=\"],\n  [\"355\u00d76=\", \"301\u00d79=\"],\n  [\"625\u00d73=\", \"923\u00d73=\"],\n  [\"839\u00d73=\", \"268\u00d76=\"],\n  [\"194\u00d77=\", \"871\u00d76=\"],\n  [\"254\u00d78=\", \"149\u00d76=\"],\n  [\"458\u00d74=\", \"984\u00d72=\"],\n  [\"410\u00d74=\", \"650\u00d74=\"],\n  [\"576\u00d78=\", \"839\u00d76=\"],\n  [\"472\u00d78=\", \"447\u00d79=\"],\n  [\"192\u00d74=\", \"445\u00d76=\"],\n  [\"205\u00d79=\", \"547\u00d77=\"],\n  [\"995\u00d79=\", \"927\u00d76=\"],\n  [\"795\u00d74=\", \"644\u00d78=\"],\n  [\"259\u00d77=\", \"904\u00d77=\"],\n  [\"890\u00d78=\", \"705\u00d74=\"],\n  [\"507\u00d79=\", \"690\u00d77=\"],\n  [\"781\u00d72=\", \"541\u00d73=\"],\n]);\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Snapshot the current text of every paragraph up front so that a\n// replacement performed on one paragraph can never influence which\n// original value we think another paragraph held.\nconst items = paragraphs.items;\nconst originalTexts = items.map((p) => p.text);\n\nfor (let i = 0; i < items.length; i++) {\n  const newText = replacements.get(originalTexts[i]);\n  if (newText !== undefined) {\n    items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every \"A\u00d7B=\" expression in the practice table.\n# Each value is updated by addressing its paragraph / table cell directly\n# (rather than via text search-and-replace), so there is no risk of a\n# later update accidentally matching text that an earlier update just\n# inserted (several \"before\" values equal other cells' \"after\" values).\n\n$d = $word.ActiveDocument\n\n# Date line is the first paragraph in the document body.\n$d.Paragraphs(1).Range.Text = \"2025-08-27 Wednesday\"\n\n$t = $d.Tables(1)\n\n# Rows 1, 5, 10, 15, 20 (1-indexed) hold the five rows of math problems;\n# the rows in between are blank spacer rows. Each row has 5 columns.\n$newValues = @{\n    1  = @(\"925\u00d73=\", \"101\u00d78=\", \"705\u00d78=\", \"205\u00d79=\", \"687\u00d75=\")\n    5  = @(\"610\u00d74=\", \"279\u00d79=\", \"160\u00d72=\", \"301\u00d79=\", \"923\u00d73=\")\n    10 = @(\"268\u00d76=\", \"871\u00d76=\", \"149\u00d76=\", \"984\u00d72=\", \"650\u00d74=\")\n    15 = @(\"839\u00d76=\", \"447\u00d79=\", \"445\u00d76=\", \"547\u00d77=\", \"927\u00d76=\")\n    20 = @(\"644\u00d78=\", \"904\u00d77=\", \"705\u00d74=\", \"690\u00d77=\", \"541\u00d73=\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $values = $newValues[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
